# Update cryptos price/volume data (scheduled GitHub Actions refresh).
# All Price (D) and Volume(1h) (E) cells are stored as text in the workbook,
# so numeric-looking Price values are entered with a leading apostrophe to
# force text interpretation, then the cell style is reset to "Normal" so
# Excel does not leave a stray text-format style applied to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.930.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.34%  "
$ws.Range("D3").Value = "'2.918.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.00%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'586.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.47%  "
$ws.Range("E6").Value = "  -5.81%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.506"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.20%  "
$ws.Range("D9").Value = "'2.916.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.04%  "
$ws.Range("D10").Value = "'6.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").Value = "'0.145"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.69%  "
$ws.Range("E12").Value = "  -4.06%  "
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.90%  "
$ws.Range("D14").Value = "'33.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.98%  "
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "'3.403.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.97%  "
$ws.Range("D17").Value = "'60.892.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.52%  "
$ws.Range("E18").Value = "  -4.48%  "
$ws.Range("D19").Value = "'2.921.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.09%  "
$ws.Range("D20").Value = "'429.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.44%  "
$ws.Range("D21").Value = "'13.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.88%  "
$ws.Range("D22").Value = "'0.682"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("E23").Value = "  -5.22%  "
$ws.Range("D24").Value = "'80.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").Value = "'2.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.20%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "'10.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.60%  "
$ws.Range("D27").Value = "'11.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.31%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "'7.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.48%  "
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("D32").Value = "'2.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").Value = "'26.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.77%  "
$ws.Range("D34").Value = "'0.107"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("E37").Value = "  -4.97%  "
$ws.Range("D38").Value = "'3.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.90%  "
$ws.Range("E39").Value = "  -2.56%  "
$ws.Range("D40").Value = "'49.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("E41").Value = "  -5.34%  "
$ws.Range("D42").Value = "'8.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.25%  "
$ws.Range("D43").Value = "'0.296"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("D44").Value = "'41.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.65%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'379.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.42%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0350"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("D47").Value = "'2.700.53"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "'132.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D50").Value = "'24.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").Value = "'0.107"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.25%  "
